# Applies the "Automatic update of files." commit:
#   - The data rows (sheet rows 2..105) are re-ordered (the source process that
#     regenerates this worksheet re-emitted the records in a different order).
#   - Column C ("Förändrad") is bumped from 46063 to 46064 for every data row.
#
# The permutation below maps each NEW row number (2..105) to the OLD row number
# that its data currently lives in, as derived from the OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# oldRowForNewRow[i] = old row number supplying the data for new row (i+2)
$oldRowForNewRow = @(2,3,4,5,6,11,12,10,7,9,8,13,14,15,16,17,22,23,30,73,92,19,32,86,65,93,81,82,58,54,20,76,63,84,57,24,66,80,50,38,103,47,105,104,74,55,101,102,61,35,95,97,98,27,31,28,26,41,87,91,96,78,34,94,36,37,42,79,43,40,51,45,46,64,48,49,59,88,21,44,68,25,70,71,85,60,99,29,52,39,77,100,72,56,75,69,67,53,90,33,89,62,18,83)

$firstDataRow = 2
$lastDataRow = 105
$numRows = $lastDataRow - $firstDataRow + 1   # 104
$numCols = 26                                  # columns A..Z

# Read the whole data block (A2:Z105) in one shot, capturing both literal
# values and formulas (HYPERLINK(...) cells) via the Formula property.
$srcRange = $ws.Range("A2:Z105")
$srcData = $srcRange.Formula

# Build the reordered block in memory.
$newData = New-Object 'object[,]' $numRows, $numCols

for ($i = 0; $i -lt $numRows; $i++) {
    $oldRow = $oldRowForNewRow[$i]
    $srcArrayRow = $oldRow - $firstDataRow + 1   # 1-based row index into $srcData
    for ($col = 1; $col -le $numCols; $col++) {
        $newData[$i, $col - 1] = $srcData[$srcArrayRow, $col]
    }
    # Column C is the 3rd column (index 2, 0-based) -> bump 46063 to 46064
    $newData[$i, 2] = 46064
}

# Write the reordered, updated block back in one shot.
$dstRange = $ws.Range("A2:Z105")
$dstRange.Formula = $newData
